$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the two
# new header cells before setting their text, so they pick up the same
# cellXfs entry (s="1") rather than minting a new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns (identical values in both
# columns for every row, rows 2-12).
$values = @(3, 5, 5, 5, 9, 9, 5, 8, 6, 4, 2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
